$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give D7 a numeric value and apply the "Normal" cell style to it. Applying
# the named style creates a brand-new cellXfs entry (color/style support),
# even though its formatting (General number format, default font/fill/border)
# mirrors the workbook's very first cell style.
$ws.Range("D7").Value = 1
$ws.Range("D7").Style = "Normal"

# Leave the active selection on E3, matching the saved view state.
$ws.Range("E3").Select()
